$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, matching the existing header style (bold,
# bordered, centered) used by the other header cells such as G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Corresponding numeric data cell for the new column
$ws.Range("H2").Value = 0
